$d = $word.ActiveDocument

# 1) Title: "Play Book of Thieves Free Slot Game | Review" -> "Play Book of Thieves for Free"
$d.Content.Find.Execute(
    "Play Book of Thieves Free Slot Game | Review", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Book of Thieves for Free", 2)

# 2) Remove the whole "Meta description: ..." paragraph entirely.
$metaFind = $d.Content
$metaFind.Find.Execute("Meta description", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
$metaPara = $d.Range($metaFind.Start, $metaFind.Start)
[void]$metaPara.Expand(4)
$metaPara.Delete()

# 3) "What we like" bullet replacements.
$d.Content.Find.Execute(
    "Well-drawn symbols and smooth animations", $true, $false, $false, $false, $false,
    $true, 1, $false, "Interesting gameplay mechanics and rules", 2)

$d.Content.Find.Execute(
    "Additional box next to the third reel creates opportunities for big wins", $true, $false, $false, $false, $false,
    $true, 1, $false, "Well-drawn symbol design and animations", 2)

$d.Content.Find.Execute(
    "Several special features and bonuses", $true, $false, $false, $false, $false,
    $true, 1, $false, "Exciting special features and bonuses", 2)

$d.Content.Find.Execute(
    "Relatively high RTP rate of 96.19%", $true, $false, $false, $false, $false,
    $true, 1, $false, "Relatively high RTP rate", 2)

# 4) "What we don't like" bullet replacements.
$d.Content.Find.Execute(
    "Average volatility may not appeal to everyone", $true, $false, $false, $false, $false,
    $true, 1, $false, "Limited betting limits", 2)

$d.Content.Find.Execute(
    "No progressive jackpot", $true, $false, $false, $false, $false,
    $true, 1, $false, "Average volatility", 2)

# 5) Insert a new bold paragraph "Play Book of Thieves for Free" right after the
#    paragraph that used to read "No progressive jackpot" (now "Average volatility").
$lastBullet = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Trim() -eq "Average volatility") {
        $lastBullet = $d.Paragraphs.Item($i)
    }
}
$lastBullet.Range.InsertParagraphAfter()
$newIndex = $lastBullet.Index + 1
$newPara = $d.Paragraphs.Item($newIndex)
$newPara.Style = "Normal"
$newPara.Range.InsertAfter("Play Book of Thieves for Free")
$textOnly = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$textOnly.Font.Bold = $true

# 6) Replace the image-prompt paragraph text.
$d.Content.Find.Execute(
    "Prompt: Design a feature image for Book of Thieves that conveys the game's medieval and mysterious atmosphere while highlighting the main character: a happy Maya warrior with glasses. Guidelines: - Create a cartoon-style image - Utilize a color scheme that fits the game's dark atmosphere - Incorporate the Maya warrior character with glasses in a way that makes them stand out - Feature the book symbol somewhere in the image to tie in with the game's theme - Consider adding elements like a medieval village or symbols from the game to enhance the image's overall design. Feel free to add your own creative flair to the image, and remember to stay true to the game's theme and tone.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Read our review of Book of Thieves and find out why you should play it for free.", 2)

Write-Output "done"
